$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "24/08/2020"
$ws.Range("B52").Value = "Data Influye"
$ws.Range("C52").Value = 1940
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 79
$ws.Range("F52").Value = 4

$ws.Range("C52:F52").NumberFormat = $ws.Range("C51:F51").NumberFormat
